$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.225.82"
$ws.Range("E2").Value = "  -1.74%  "
$ws.Range("D3").Value = "2.911.24"
$ws.Range("E3").Value = "  -3.02%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "376.94"
$ws.Range("E5").Value = "  +6.50%  "
$ws.Range("D6").Value = "102.31"
$ws.Range("E6").Value = "  -4.58%  "
$ws.Range("E7").Value = "  -3.15%  "
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "0.584"
$ws.Range("E9").Value = "  -4.24%  "
$ws.Range("D10").Value = "36.79"
$ws.Range("E10").Value = "  -3.66%  "
$ws.Range("E11").Value = "  -0.23%  "
$ws.Range("D12").Value = "0.0831"
$ws.Range("E12").Value = "  -2.84%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "18.18"
$ws.Range("E13").Value = "  -4.67%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "3.370.94"
$ws.Range("E14").Value = "  -2.85%  "
$ws.Range("D15").Value = "7.30"
$ws.Range("E15").Value = "  -4.05%  "
$ws.Range("D16").Value = "2.913.20"
$ws.Range("E16").Value = "  -2.61%  "
$ws.Range("D17").Value = "0.924"
$ws.Range("E17").Value = "  -9.34%  "
$ws.Range("D18").Value = "51.158.30"
$ws.Range("E18").Value = "  -1.91%  "
$ws.Range("D19").Value = "3.40"
$ws.Range("E19").Value = "  -0.34%  "
$ws.Range("D20").Value = "7.31"
$ws.Range("E20").Value = "  -2.47%  "
$ws.Range("D21").Value = "12.84"
$ws.Range("E21").Value = "  -5.25%  "
$ws.Range("D22").Value = "0.0₃0941"
$ws.Range("E22").Value = "  -3.04%  "
$ws.Range("D23").Value = "68.06"
$ws.Range("E23").Value = "  -1.66%  "
$ws.Range("D24").Value = "259.99"
$ws.Range("E24").Value = "  -1.74%  "
$ws.Range("D25").Value = "2.75"
$ws.Range("E25").Value = "  +0.95%  "
$ws.Range("D26").Value = "0.168"
$ws.Range("E26").Value = "  -4.67%  "
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("B28").Value = "LEO"
$ws.Range("C28").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D28").Value = "4.11"
$ws.Range("E28").Value = "  -5.05%  "
$ws.Range("D29").Value = "25.56"
$ws.Range("E29").Value = "  -4.91%  "
$ws.Range("E30").Value = "  -4.30%  "
$ws.Range("D31").Value = "6.77"
$ws.Range("E31").Value = "  +4.75%  "
$ws.Range("E32").Value = "  -4.41%  "
$ws.Range("D33").Value = "9.74"
$ws.Range("E33").Value = "  -4.84%  "
$ws.Range("D34").Value = "2.10"
$ws.Range("E34").Value = "  -3.50%  "
$ws.Range("D35").Value = "51.25"
$ws.Range("E35").Value = "  +0.33%  "
$ws.Range("D36").Value = "33.81"
$ws.Range("E36").Value = "  -6.21%  "
$ws.Range("E37").Value = "  +0.34%  "
$ws.Range("D38").Value = "0.0417"
$ws.Range("E38").Value = "  -4.51%  "
$ws.Range("D39").Value = "2.98"
$ws.Range("E39").Value = "  -10.10%  "
$ws.Range("D40").Value = "16.83"
$ws.Range("E40").Value = "  -4.38%  "
$ws.Range("D41").Value = "2.53"
$ws.Range("E41").Value = "  -9.18%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").Value = "0.114"
$ws.Range("E42").Value = "  -2.80%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "1.80"
$ws.Range("E43").Value = "  -9.09%  "
$ws.Range("D44").Value = "122.54"
$ws.Range("E44").Value = "  -1.33%  "
$ws.Range("D45").Value = "21.42"
$ws.Range("E45").Value = "  -5.86%  "
$ws.Range("E46").Value = "  -3.33%  "
$ws.Range("D47").Value = "0.268"
$ws.Range("E47").Value = "  +10.49%  "
$ws.Range("D48").Value = "2.023.92"
$ws.Range("E48").Value = "  -4.63%  "
$ws.Range("E49").Value = "  -2.73%  "
$ws.Range("D50").Value = "3.13"
$ws.Range("E50").Value = "  -5.96%  "
$ws.Range("D51").Value = "3.204.76"
$ws.Range("E51").Value = "  -2.66%  "
